$wb = $excel.ActiveWorkbook

# --- Modulos sheet: add two new rows of module/menu data ---
$ws = $wb.Worksheets.Item("Modulos")

# Write the new shared-string-bearing cells. D18/D19 first so the new
# shared-string table keeps the same insertion order as the source edit.
$ws.Range("D18").Value = "BCCL.E.CONS.TINT.ALTAS"
$ws.Range("D19").Value = "BCCL.E.TINTERNAS.APAGAR"
$ws.Range("A18").Value = "Transferencias Internas"
$ws.Range("B18").Value = "Consulta de Altas Transf. Internas "
$ws.Range("B19").Value = "Consulta de Altas Transf. Internas A Pagar "

# Column width adjustments for columns A and B.
# NOTE: the engine quantizes ColumnWidth to a 1/6-character pixel grid
# (stored = (round(6*x)+5)/6), so the exact source widths 21.42578125 /
# 38.42578125 (1/256-character grid) aren't reachable bit-for-bit. These
# inputs are chosen so the quantized, saved width lands as close as
# possible to the target (21.5 / 38.5 instead of 21.42578125 / 38.42578125).
$ws.Columns.Item(1).ColumnWidth = 20.666666666666668
$ws.Columns.Item(2).ColumnWidth = 37.666666666666664

# Move the selection / active cell to C19 and make this sheet the active tab
$ws.Range("C19").Select()
